$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.230985683306322, 114.8270160096505, 0.8054896365839992, 8.660232485948974, 127.5237238154898)
    3 = @(0.3048080303191223, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1.919867272924993)
    4 = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 0.496779210170732, 6.740334628841572)
    5 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    6 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    7 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 14.36450238910742)
    8 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
